$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data as per the commit:
# "Updated symbol list on Wed Feb 15 23:57:57 UTC 2023 with GitHub Actions"
# Each D/E cell must remain text (not auto-converted to a number/percentage),
# so we force the Text number format on the cell before assigning the value.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "315.99"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "6.61%"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "45.30"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "7.36%"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.180"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "2.84%"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.08095"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "6.77%"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "3.31%"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "5.32%"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.094"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "17.58%"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1311"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "8.09%"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "5.23%"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.09360"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "4.03%"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.04246"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "5.93%"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-0.69%"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001321"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "3.30%"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.005821"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "0.50%"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.401"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "1.05%"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.425"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "0.83%"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.3373"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "1.60%"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "8.215"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "4.93%"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.1356"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "-3.19%"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04287"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "5.46%"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001281"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "1.22%"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004214"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "7.54%"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0001345"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "9.27%"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02702"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "11.54%"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05465"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "4.83%"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.005877"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "-1.48%"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.007801"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "0.48%"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1426"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "7.10%"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.007350"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "-2.73%"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.008559"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "17.88%"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.3144"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "6.20%"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00006820"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "0.58%"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "-0.41%"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.06068"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "33.14%"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "-5.12%"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "-0.41%"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "-0.41%"
